$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.603.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("E6").Value = '  -3.01%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '33.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.98%  '
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.054.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.808.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.642'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.506.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '256.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0751'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.15%  '
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.01%  '
$ws.Range("E28").Value = '  -2.88%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0516'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("E34").Value = '  +6.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.460.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.19%  '
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.634'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("E39").Value = '  +2.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.899'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0508'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.953.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("E48").Value = '  +0.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.998'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.65%  '
